$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("H2").Value = 1

# Row 3
$ws.Range("F3").Value = "DİYARBAKIR"
$ws.Range("H3").Value = 1

# Row 4
$ws.Range("F4").Value = "KARABÜK"

# Row 5
$ws.Range("F5").Value = "GİRESUN"

# Row 6
$ws.Range("F6").Value = "İZMİR"
$ws.Range("H6").Value = 2

# Row 7
$ws.Range("F7").Value = "GİRESUN"
$ws.Range("H7").Value = 3

# Row 8
$ws.Range("F8").Value = "KARABÜK"
$ws.Range("H8").Value = 1

# Row 9
$ws.Range("F9").Value = "DİYARBAKIR"

# Row 11
$ws.Range("F11").Value = "İSTANBUL"
